# Restore original "Recorded By" ordering for selected rows in the attendance report.
# Each listed G-cell's comma-separated recorder list is reordered to match the
# upstream canonical value, while all other G-cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'System, backup@backdoor.com, system'
$ws.Range('G4').Value = 'System, backup@backdoor.com'
$ws.Range('G5').Value = 'System, backup@backdoor.com'
$ws.Range('G7').Value = 'admin@admin.com, System'
$ws.Range('G8').Value = 'System, backup@backdoor.com'
$ws.Range('G11').Value = 'dnasr281@gmail.com, System'
$ws.Range('G17').Value = 'dnasr281@gmail.com, System'
$ws.Range('G28').Value = 'System, backup@backdoor.com, system'
$ws.Range('G30').Value = 'System, backup@backdoor.com'
$ws.Range('G31').Value = 'System, backup@backdoor.com'
$ws.Range('G33').Value = 'admin@admin.com, System'
$ws.Range('G34').Value = 'System, backup@backdoor.com'
$ws.Range('G37').Value = 'dnasr281@gmail.com, System'
$ws.Range('G43').Value = 'dnasr281@gmail.com, System'
$ws.Range('G54').Value = 'System, backup@backdoor.com, system'
$ws.Range('G56').Value = 'System, backup@backdoor.com'
$ws.Range('G57').Value = 'System, backup@backdoor.com'
$ws.Range('G59').Value = 'admin@admin.com, System'
$ws.Range('G60').Value = 'System, backup@backdoor.com'
$ws.Range('G63').Value = 'dnasr281@gmail.com, System'
$ws.Range('G69').Value = 'dnasr281@gmail.com, System'
$ws.Range('G80').Value = 'System, backup@backdoor.com'
$ws.Range('G81').Value = 'System, backup@backdoor.com'
$ws.Range('G82').Value = 'System, backup@backdoor.com'
$ws.Range('G93').Value = 'dnasr281@gmail.com, System'
$ws.Range('G94').Value = 'dnasr281@gmail.com, System'
$ws.Range('G96').Value = 'dnasr281@gmail.com, System'
$ws.Range('G106').Value = 'System, backup@backdoor.com'
$ws.Range('G107').Value = 'System, backup@backdoor.com'
$ws.Range('G108').Value = 'System, backup@backdoor.com'
$ws.Range('G119').Value = 'dnasr281@gmail.com, System'
$ws.Range('G120').Value = 'dnasr281@gmail.com, System'
$ws.Range('G122').Value = 'dnasr281@gmail.com, System'
$ws.Range('G132').Value = 'System, backup@backdoor.com'
$ws.Range('G133').Value = 'System, backup@backdoor.com'
$ws.Range('G134').Value = 'System, backup@backdoor.com'
$ws.Range('G145').Value = 'dnasr281@gmail.com, System'
$ws.Range('G146').Value = 'dnasr281@gmail.com, System'
$ws.Range('G148').Value = 'dnasr281@gmail.com, System'
